$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column "ano" header at A2
$ws.Range("A2").Value = "ano"

# Fill A3:A14 with year 2023
for ($r = 3; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = 2023
}

# Update selection to match the diff (A4:A14 active cell A4)
$ws.Range("A4:A14").Select()
